$d = $word.ActiveDocument

# --- Change 1: Title text - remove "Assessment" (trailing space kept) ---
$d.Content.Find.Execute(
    "Security Report: OWASP Top 10 Security Risks Assessment",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Security Report: OWASP Top 10 Security Risks ", 2
) | Out-Null

# --- Change 2: "make it safe" -> "make it more safe" (as 3 runs) ---
$findRng = $d.Content.Duplicate
$findRng.Find.Execute("make it safe", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertAt = $findRng.Start + 8
$insPoint = $d.Range($insertAt, $insertAt)
$insPoint.InsertAfter("more ") | Out-Null
$newWordRng = $d.Range($insertAt, $insertAt + 5)
$newWordRng.Font.Size = 99
$newWordRng.Font.SizeBi = 99
$newWordRng.Font.Size = 12
$newWordRng.Font.SizeBi = 12

# --- Change 3: "References" Heading1 paragraph gets sz=30/szCs=30 ---
$refHeading = $d.Paragraphs.Item(84)
$refHeadingRng = $refHeading.Range
$refHeadingRng.Font.Size = 15
$refHeadingRng.Font.SizeBi = 15

# --- Change 4: Hyperlink paragraph gets sz=24/szCs=24 on pPr and each run ---
$refLink = $d.Paragraphs.Item(85)
$refLinkRng = $refLink.Range
$refLinkRng.Font.Size = 12
$refLinkRng.Font.SizeBi = 12

Write-Output "All changes applied"
